# selectWindow.xlsx update -> v0.4
# Rebuilds the TestCase sheet: header row gains storeText/echo command columns,
# the click step now carries its target as JSON ({"target":"id=btn1"}), and the
# single "click" row is expanded into three rows (btn1/btn2/btn3) followed by a
# storeText + echo round trip and a tab/title/url based selectWindow example.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header row (row 1): C1/D1 keep their existing text+style, E1/F1 are new
#    "storeText"/"echo" headers that reuse the existing header look (border +
#    center/middle alignment) but in a slightly different font, matching the
#    other command-name headers' family.
# ---------------------------------------------------------------------------
$ws.Range("E1").Copy() | Out-Null
$ws.Range("E1").PasteSpecial(-4122) | Out-Null       # xlPasteFormats
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null       # xlPasteFormats
$ws.Range("E1:F1").Font.Name = "ＭＳ Ｐゴシック"
$ws.Range("E1").Value = "storeText"
$ws.Range("F1").Value = "echo"

# ---------------------------------------------------------------------------
# 2. Row 2: C2 loses its old "id=btn1" value (now lives in the data rows as
#    JSON), D2 picks up the same (unused) hyperlink-flavoured font that B2
#    uses, and E2/F2 get the new storeText argument + echo argument.
# ---------------------------------------------------------------------------
$ws.Range("C2").Value = ""

$ws.Range("B2").Copy() | Out-Null
$ws.Range("D2").PasteSpecial(-4122) | Out-Null       # xlPasteFormats

$ws.Range("A2").Copy() | Out-Null
$ws.Range("F2").PasteSpecial(-4122) | Out-Null       # xlPasteFormats

$ws.Range("E2").Value = '{"target":"id=result","value":"result"}'
$ws.Range("F2").Value = '${result}'

# ---------------------------------------------------------------------------
# 3. Row 3 becomes the first of three "click" rows; E3 loses the stray 2000
#    and gains an F3 companion cell (matching E3's style).
# ---------------------------------------------------------------------------
$ws.Range("E3").Value = ""

$ws.Range("E3").Copy() | Out-Null
$ws.Range("F3").PasteSpecial(-4122) | Out-Null       # xlPasteFormats

$ws.Range("C3").Value = '{"target":"id=btn1"}'
$ws.Range("D3").Value = '{"target":"tab=1"}'

# ---------------------------------------------------------------------------
# 4. Rows 4 and 5 are brand-new rows, cloned (format-only) from row 3 so they
#    pick up the exact same borders/alignment/fonts, then filled with data.
# ---------------------------------------------------------------------------
$ws.Range("A3:F3").Copy() | Out-Null
$ws.Range("A4").PasteSpecial(-4122) | Out-Null       # xlPasteFormats
$ws.Range("A3:F3").Copy() | Out-Null
$ws.Range("A5").PasteSpecial(-4122) | Out-Null       # xlPasteFormats

$ws.Range("A4").Value = 2
$ws.Range("C4").Value = '{"target":"id=btn2"}'
$ws.Range("D4").Value = '{"target":"byTitle"}'

$ws.Range("A5").Value = 3
$ws.Range("C5").Value = '{"target":"id=btn3"}'
$ws.Range("D5").Value = '{"target":"http://localhost:9001/selectWindow/byUrl.html"}'

# ---------------------------------------------------------------------------
# 5. Column widths: C/D grow to fit their longer JSON targets, E/F are new.
# ---------------------------------------------------------------------------
$ws.Columns("C").ColumnWidth = 18.429
$ws.Columns("D").ColumnWidth = 56
$ws.Columns("E").ColumnWidth = 34.429
$ws.Columns("F").ColumnWidth = 8.571

# ---------------------------------------------------------------------------
# 6. Selection: leave the cursor where Excel parked it after the edits.
# ---------------------------------------------------------------------------
$ws.Range("D10").Select() | Out-Null
